# Add 2022-Q3 data
# 1) Insert a new "2022-Q3" worksheet right after "总计", built from a copy
#    of "2022-Q2" so it inherits matching formatting (header/index-column
#    styles, column count, page margins, etc.), then overwrite its values.
# 2) Prepend a new row to the "总计" summary sheet for the 2022-Q3 totals,
#    shifting the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# --- Step 1: create the new "2022-Q3" sheet -------------------------------
$q2Sheet.Copy($null, $totalSheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# The copied sheet has 8 data rows (rows 2-9); 2022-Q3 only has 4 (rows 2-5).
$newSheet.Rows("6:9").Delete()

# Columns that hold text-like values (fund code with leading zeros, and
# decimal-looking numbers that are actually stored as text) must be forced
# to Text format so the leading zeros / exact string values survive.
$newSheet.Range("B2:B5").NumberFormat = "@"
$newSheet.Range("D2:G5").NumberFormat = "@"

$q3Data = @(
    @("014179", "中银证券远见价值混合A", "1.56", "93.65", "5.37", "0.0838", 4),
    @("003980", "中银证券瑞益灵活配置混合A", "0.66", "91.21", "3.68", "0.0243", 6),
    @("014180", "中银证券远见价值混合C", "0.16", "93.65", "5.37", "0.0086", 4),
    @("003981", "中银证券瑞益灵活配置混合C", "0.19", "91.21", "3.68", "0.0070", 6)
)

for ($i = 0; $i -lt $q3Data.Length; $i++) {
    $r = $i + 2
    $row = $q3Data[$i]
    $newSheet.Cells.Item($r, 1).Value = $i
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $newSheet.Cells.Item($r, 5).Value = $row[3]
    $newSheet.Cells.Item($r, 6).Value = $row[4]
    $newSheet.Cells.Item($r, 7).Value = $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

# --- Step 2: add the new row to "总计" -------------------------------------
$totalData = @(
    @(0, "2022-Q3", 4, 0.12),
    @(1, "2022-Q2", 8, 1.64),
    @(2, "2022-Q1", 10, 1.77),
    @(3, "2021-Q4", 7, 1.28),
    @(4, "2021-Q3", 8, 4.68),
    @(5, "2021-Q2", 8, 6.2),
    @(6, "2021-Q1", 5, 6.51),
    @(7, "2020-Q4", 6, 12.66)
)

for ($i = 0; $i -lt $totalData.Length; $i++) {
    $r = $i + 2
    $row = $totalData[$i]
    $totalSheet.Cells.Item($r, 1).Value = $row[0]
    $totalSheet.Cells.Item($r, 2).Value = $row[1]
    $totalSheet.Cells.Item($r, 3).Value = $row[2]
    $totalSheet.Cells.Item($r, 4).Value = $row[3]
}

# Row 9 of "总计" is brand new; copy column-A formatting from row 8 so the
# index cell keeps the bold/bordered "s=2" style used throughout the column.
$totalSheet.Cells.Item(8, 1).Copy()
$totalSheet.Cells.Item(9, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
